$wb = $excel.ActiveWorkbook

# Status for the two handoff rows moved from "Ready for handoff" to
# "In Translation" on every sheet (Overview's zh-cn/de-de status columns,
# plus the Status column on each per-locale report sheet).
foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($name -eq "Overview") {
        $cols = @("E", "F")
    } else {
        $cols = @("C")
    }
    foreach ($col in $cols) {
        foreach ($row in 2..3) {
            $cell = $ws.Range("$col$row")
            if ($cell.Value() -eq "Ready for handoff") {
                $cell.Value = "In Translation"
            }
        }
    }
}

# The status columns narrowed once the new, shorter status text was in
# place (re-fit to contents). Apply the narrower width to those same
# columns on every sheet.
foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($name -eq "Overview") {
        $cols = @(5, 6)
    } else {
        $cols = @(3)
    }
    foreach ($colIndex in $cols) {
        $ws.Columns.Item($colIndex).ColumnWidth = 12.5
    }
}
